$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "NB" row (row 8). This shifts the old SVM row (9) up to row 8,
# carrying over its existing style (bold/bordered column A, etc).
$ws.Rows.Item(8).Delete()

# Extend header formatting (bold/border) from an existing styled header cell
# (B1) into the three new header columns F1:H1.
$ws.Range("B1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# ---- Header row ----
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# ---- Row 2 : LR ----
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.9188231243875957
$ws.Range("D2").Value = 0.01138684070920371
$ws.Range("E2").Value = 0.9129633749106221
$ws.Range("F2").Value = 0.008360109853643438
$ws.Range("G2").Value = 0.9061351659119197
$ws.Range("H2").Value = 0.008709544003751301

# ---- Row 3 : LDA ----
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.9152453589682477
$ws.Range("D3").Value = 0.009912140208083203
$ws.Range("E3").Value = 0.9001120203384444
$ws.Range("F3").Value = 0.007484977942516059
$ws.Range("G3").Value = 0.9023902968671381
$ws.Range("H3").Value = 0.009633082692969813

# ---- Row 4 : KNN ----
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.8778281825163529
$ws.Range("D4").Value = 0.02026078736925795
$ws.Range("E4").Value = 0.8929525171473214
$ws.Range("F4").Value = 0.009900576506895261
$ws.Range("G4").Value = 0.8757122427901803
$ws.Range("H4").Value = 0.009343281117336819

# ---- Row 5 : DTREE (was CART) ----
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.9134538280236223
$ws.Range("D5").Value = 0.0101221123050728
$ws.Range("E5").Value = 0.9274420698604381
$ws.Range("F5").Value = 0.008574366708933129
$ws.Range("G5").Value = 0.9041818278117633
$ws.Range("H5").Value = 0.007838049613821689

# ---- Row 6 : RTREE ----
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8849858319430099
$ws.Range("D6").Value = 0.01663871656508185
$ws.Range("E6").Value = 0.8395913773469982
$ws.Range("F6").Value = 0.01216175278667438
$ws.Range("G6").Value = 0.8683909324435264
$ws.Range("H6").Value = 0.008891532393564765

# ---- Row 7 : XTREE ----
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.9328118958714017
$ws.Range("D7").Value = 0.01293349059205216
$ws.Range("E7").Value = 0.9245133868276794
$ws.Range("F7").Value = 0.01071376677207109
$ws.Range("G7").Value = 0.9168703159344297
$ws.Range("H7").Value = 0.008647476777694138

# ---- Row 8 : SVM ----
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.9217515425968592
$ws.Range("D8").Value = 0.009948321308492932
$ws.Range("E8").Value = 0.9300450199941738
$ws.Range("F8").Value = 0.004911441279662837
$ws.Range("G8").Value = 0.9155681788088239
$ws.Range("H8").Value = 0.009247709485030942
